$wb = $excel.ActiveWorkbook

# The status text for the 725f5bfa... row changed from "Ready for handoff"
# to "Handback transform failed". That text is shared by every sheet that
# shows this row's status: Overview!B3 (zh-cn column), Overview!C3
# (de-de column), and the "Status" column (C3) on both the zh-cn and
# de-de detail sheets.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("B3").Value = "Handback transform failed"
$wsOverview.Range("C3").Value = "Handback transform failed"
$wsZhCn.Range("C3").Value = "Handback transform failed"
$wsDeDe.Range("C3").Value = "Handback transform failed"

# Sheet "zh-cn": add an Error Detail value in column L, row 3 (the
# 725f5bfa... row) explaining the handback/handoff filename mismatch.
$wsZhCn.Range("L3").Value = "Handback file name: 4aial15c.1n2 is different with handoff file name: 725f5bfa-059c-4bce-a02f-3523664e5985.5d8bf74cb3651a637180bb72d451fa858ea9a03e.zh-cn."

# Sheet "de-de": same Error Detail addition for its own locale.
$wsDeDe.Range("L3").Value = "Handback file name: 4aial15c.1n2 is different with handoff file name: 725f5bfa-059c-4bce-a02f-3523664e5985.5d8bf74cb3651a637180bb72d451fa858ea9a03e.de-de."
